# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as described by the authoritative diff (columns H-N, various rows).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 2 (G2=5489)
$ws.Range("H2").Value = 110.25
$ws.Range("I2").Value = 123.5
$ws.Range("K2").Value = 123.5
$ws.Range("M2").Value = -10.5
# row 5 (G5=5503)
$ws.Range("H5").Value = 102.92857
$ws.Range("I5").Value = 94.63636
$ws.Range("J5").Value = 133.33333
$ws.Range("K5").Value = 94.63636
$ws.Range("L5").Value = 133.33333
$ws.Range("M5").Value = 20.36364
$ws.Range("N5").Value = -363.33333
# row 95 (G95=18200)
$ws.Range("H95").Value = 24664
$ws.Range("J95").Value = 24664
$ws.Range("L95").Value = 24664
$ws.Range("N95").Value = -30156
# row 116 (G116=27778)
$ws.Range("H116").Value = 3286
$ws.Range("I116").Value = 3009.5715
$ws.Range("J116").Value = 3673
$ws.Range("K116").Value = 3009.5715
$ws.Range("L116").Value = 3673
$ws.Range("M116").Value = 432.4285
$ws.Range("N116").Value = -10557
# row 129 (G129=36115)
$ws.Range("H129").Value = 997.8570999999999
$ws.Range("I129").Value = 452.63635
$ws.Range("J129").Value = 2997
$ws.Range("K129").Value = 1357.90905
$ws.Range("L129").Value = 8991
$ws.Range("M129").Value = 3642.09095
$ws.Range("N129").Value = -18991
# row 137 (G137=44013)
$ws.Range("H137").Value = 2410.818
$ws.Range("I137").Value = 1540.8334
$ws.Range("J137").Value = 2737.0625
$ws.Range("K137").Value = 4622.5002
$ws.Range("L137").Value = 8211.1875
$ws.Range("M137").Value = -2072.5002
$ws.Range("N137").Value = -13311.1875
# row 138 (G138=44169)
$ws.Range("H138").Value = 3742
$ws.Range("I138").Value = 1023.5
$ws.Range("J138").Value = 5554.3335
$ws.Range("K138").Value = 3070.5
$ws.Range("L138").Value = 16663.0005
$ws.Range("M138").Value = 2069.5
$ws.Range("N138").Value = -26943.0005

$ws = $wb.Worksheets.Item("ARM")
# row 6 (G6=2226)
$ws.Range("H6").Value = 14288714
$ws.Range("I6").Value = 12004200
$ws.Range("J6").Value = 19999998
$ws.Range("K6").Value = 12004200
$ws.Range("L6").Value = 19999998
$ws.Range("M6").Value = -12004027
$ws.Range("N6").Value = -20000344
# row 124 (G124=34252)
$ws.Range("H124").Value = 16809
$ws.Range("J124").Value = 16809
$ws.Range("L124").Value = 16809
$ws.Range("N124").Value = -26629
# row 132 (G132=43997)
$ws.Range("H132").Value = 1423
$ws.Range("I132").Value = 1423
$ws.Range("K132").Value = 4269
$ws.Range("M132").Value = -1739

$ws = $wb.Worksheets.Item("BSM")
# row 20 (G20=14149)
$ws.Range("H20").Value = 6254.25
$ws.Range("I20").Value = 6999.5
$ws.Range("J20").Value = 5509
$ws.Range("K20").Value = 6999.5
$ws.Range("L20").Value = 5509
$ws.Range("M20").Value = -6752.5
$ws.Range("N20").Value = -6003
# row 88 (G88=10626)
$ws.Range("H88").Value = 27943.6
$ws.Range("J88").Value = 27943.6
$ws.Range("L88").Value = 27943.6
$ws.Range("N88").Value = -28755.6
# row 91 (G91=10626)
$ws.Range("H91").Value = 27943.6
$ws.Range("J91").Value = 27943.6
$ws.Range("L91").Value = 27943.6
$ws.Range("N91").Value = -30751.6

$ws = $wb.Worksheets.Item("CRP")
# row 7 (G7=5361)
$ws.Range("H7").Value = 61.173912
$ws.Range("I7").Value = 44.2
$ws.Range("K7").Value = 44.2
$ws.Range("M7").Value = 68.8
# row 55 (G55=1855)
$ws.Range("H55").Value = 15356.333
$ws.Range("I55").Value = 7073
$ws.Range("J55").Value = 19498
$ws.Range("K55").Value = 7073
$ws.Range("L55").Value = 19498
$ws.Range("M55").Value = -6758
$ws.Range("N55").Value = -20128
# row 62 (G62=12580)
$ws.Range("H62").Value = 2496
$ws.Range("J62").Value = 2496
$ws.Range("L62").Value = 2496
$ws.Range("N62").Value = -3744
# row 65 (G65=12580)
$ws.Range("H65").Value = 2496
$ws.Range("J65").Value = 2496
$ws.Range("L65").Value = 12480
$ws.Range("N65").Value = -18720
# row 123 (G123=35334)
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null

$ws = $wb.Worksheets.Item("CUL")
# row 7 (G7=4728)
$ws.Range("H7").Value = 29.5
$ws.Range("I7").Value = 7.2
$ws.Range("J7").Value = 66.666664
$ws.Range("K7").Value = 21.6
$ws.Range("L7").Value = 199.999992
$ws.Range("M7").Value = 90.40000000000001
$ws.Range("N7").Value = -423.999992
# row 34 (G34=4749)
$ws.Range("H34").Value = 1953.9333
$ws.Range("I34").Value = 249
$ws.Range("J34").Value = 2075.7144
$ws.Range("K34").Value = 747
$ws.Range("L34").Value = 6227.1432
$ws.Range("M34").Value = -663
$ws.Range("N34").Value = -6395.1432
# row 36 (G36=4732)
$ws.Range("H36").Value = 5115.6665
$ws.Range("I36").Value = 173.5
$ws.Range("J36").Value = 15000
$ws.Range("K36").Value = 520.5
$ws.Range("L36").Value = 45000
$ws.Range("M36").Value = -351.5
$ws.Range("N36").Value = -45338
# row 44 (G44=4702)
$ws.Range("H44").Value = 3317.818
$ws.Range("I44").Value = 214
$ws.Range("K44").Value = 642
$ws.Range("M44").Value = -244
# row 47 (G47=4663)
$ws.Range("H47").Value = 4551.778
$ws.Range("I47").Value = 161
$ws.Range("K47").Value = 483
$ws.Range("M47").Value = -52
# row 113 (G113=27843)
$ws.Range("H113").Value = 1552.421
$ws.Range("I113").Value = 1017.5714
$ws.Range("J113").Value = 1864.4166
$ws.Range("K113").Value = 3052.7142
$ws.Range("L113").Value = 5593.2498
$ws.Range("M113").Value = -882.7142000000003
$ws.Range("N113").Value = -9933.2498
# row 139 (G139=44102)
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").Value = $null

$ws = $wb.Worksheets.Item("GSM")
# row 2 (G2=5062)
$ws.Range("H2").Value = 258.68182
$ws.Range("I2").Value = 199.44444
$ws.Range("J2").Value = 299.69232
$ws.Range("K2").Value = 199.44444
$ws.Range("L2").Value = 299.69232
$ws.Range("M2").Value = -86.44443999999999
$ws.Range("N2").Value = -525.69232
# row 11 (G11=4422)
$ws.Range("H11").Value = 12314054
$ws.Range("I11").Value = 7502052
$ws.Range("J11").Value = 26750062
$ws.Range("K11").Value = 7502052
$ws.Range("L11").Value = 26750062
$ws.Range("M11").Value = -7501913
$ws.Range("N11").Value = -26750340

$ws = $wb.Worksheets.Item("LTW")
# row 22 (G22=5277)
$ws.Range("H22").Value = 1345.5555
$ws.Range("I22").Value = 722
$ws.Range("J22").Value = 2125
$ws.Range("K22").Value = 722
$ws.Range("L22").Value = 2125
$ws.Range("M22").Value = -427
$ws.Range("N22").Value = -2715
# row 27 (G27=5277)
$ws.Range("H27").Value = 1345.5555
$ws.Range("I27").Value = 722
$ws.Range("J27").Value = 2125
$ws.Range("K27").Value = 722
$ws.Range("L27").Value = 2125
$ws.Range("M27").Value = -615
$ws.Range("N27").Value = -2339
# row 35 (G35=1697)
$ws.Range("H35").Value = 1489.2222
$ws.Range("I35").Value = 1287.875
$ws.Range("J35").Value = 3100
$ws.Range("K35").Value = 1287.875
$ws.Range("L35").Value = 3100
$ws.Range("M35").Value = -951.875
$ws.Range("N35").Value = -3772

$ws = $wb.Worksheets.Item("WVR")
# row 6 (G6=3000)
$ws.Range("H6").Value = 566.6667
$ws.Range("I6").Value = 550
$ws.Range("K6").Value = 550
$ws.Range("M6").Value = -435
# row 11 (G11=3001)
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").Value = $null
# row 18 (G18=3543)
$ws.Range("H18").Value = 10666.167
$ws.Range("I18").Value = 9999
$ws.Range("K18").Value = 9999
$ws.Range("M18").Value = -9826
# row 58 (G58=3187)
$ws.Range("H58").Value = 3128.75
$ws.Range("I58").Value = 3128.75
$ws.Range("K58").Value = 3128.75
$ws.Range("M58").Value = -2820.75
# row 107 (G107=27746)
$ws.Range("H107").Value = 41667452
$ws.Range("I107").Value = 41667452
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 125002356
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -125000436
$ws.Range("N107").Value = $null
# row 136 (G136=44031)
$ws.Range("H136").Value = 2824.2273
$ws.Range("I136").Value = 2045.2858
$ws.Range("K136").Value = 6135.857400000001
$ws.Range("M136").Value = -3585.857400000001
